$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.224441849520701
$ws.Range("C2").Value = 0.4422168556887982
$ws.Range("D2").Value = 0.07928524414171534
$ws.Range("E2").Value = 0.420912291607948
$ws.Range("G2").Value = 0.8331245288309219
$ws.Range("H2").Value = 0.7781348121798999
$ws.Range("I2").Value = 0.5395731323089805

$ws.Range("B3").Value = 1.080384171106516
$ws.Range("C3").Value = 0.3854780840128456
$ws.Range("D3").Value = 0.0717832443492199
$ws.Range("E3").Value = 0.3668040333107143
$ws.Range("G3").Value = 0.7912173271899405
$ws.Range("H3").Value = 0.7658778946088773
$ws.Range("I3").Value = 0.5336542410614769

$ws.Range("B4").Value = 0.9920831924560503
$ws.Range("C4").Value = 0.3506509177374255
$ws.Range("D4").Value = 0.06721998986452604
$ws.Range("E4").Value = 0.3337230811786895
$ws.Range("G4").Value = 0.7663798871778909
$ws.Range("H4").Value = 0.759035056955355
$ws.Range("I4").Value = 0.530590331385568

$ws.Range("B5").Value = 0.9561366869488666
$ws.Range("C5").Value = 0.3364601464040788
$ws.Range("D5").Value = 0.06537102845226173
$ws.Range("E5").Value = 0.3202748146432839
$ws.Range("G5").Value = 0.7564792546534136
$ws.Range("H5").Value = 0.7564165988064531
$ws.Range("I5").Value = 0.5294835121654629

$ws.Range("B6").Value = 0.9501699923242199
$ws.Range("C6").Value = 0.3341038435743826
$ws.Range("D6").Value = 0.065064644020012
$ws.Range("E6").Value = 0.3180436112658214
$ws.Range("G6").Value = 0.7548484855350353
$ws.Range("H6").Value = 0.7559920277196284
$ws.Range("I6").Value = 0.5293082385404375

$ws.Range("B7").Value = 0.9915982574891586
$ws.Range("C7").Value = 0.3504595310654395
$ws.Range("D7").Value = 0.06719501147472329
$ws.Range("E7").Value = 0.3335415857384447
$ws.Range("G7").Value = 0.7662454748208347
$ws.Range("H7").Value = 0.7589990571801764
$ws.Range("I7").Value = 0.5305748326840245

$ws.Range("B8").Value = 1.17473868244997
$ws.Range("C8").Value = 0.4226504255448162
$ws.Range("D8").Value = 0.07668949851232298
$ws.Range("E8").Value = 0.402224423252477
$ws.Range("G8").Value = 0.8184870500026022
$ws.Range("H8").Value = 0.7737658568335348
$ws.Range("I8").Value = 0.5374129244749994

$ws.Range("B9").Value = 1.535146559177406
$ws.Range("C9").Value = 0.5643588085901001
$ws.Range("D9").Value = 0.09566000308355171
$ws.Range("E9").Value = 0.5381854486696369
$ws.Range("G9").Value = 0.9282071144384645
$ws.Range("H9").Value = 0.8082212371615753
$ws.Range("I9").Value = 0.5554240937004522

$ws.Range("B10").Value = 1.800843055259918
$ws.Range("C10").Value = 0.6686499495361886
$ws.Range("D10").Value = 0.1098284020873592
$ws.Range("E10").Value = 0.639070785498788
$ws.Range("G10").Value = 1.013515835252917
$ws.Range("H10").Value = 0.837000441137576
$ws.Range("I10").Value = 0.5715704776142445

$ws.Range("B11").Value = 1.921943027878456
$ws.Range("C11").Value = 0.7161542255649351
$ws.Range("D11").Value = 0.1163277447645044
$ws.Range("E11").Value = 0.6852306832343089
$ws.Range("G11").Value = 1.053402763484911
$ws.Range("H11").Value = 0.8508701930509801
$ws.Range("I11").Value = 0.5795722499284039

$ws.Range("B12").Value = 1.967836055889393
$ws.Range("C12").Value = 0.7341532637036607
$ws.Range("D12").Value = 0.1187969247149709
$ws.Range("E12").Value = 0.7027524360483852
$ws.Range("G12").Value = 1.06866674214001
$ws.Range("H12").Value = 0.8562361722395053
$ws.Range("I12").Value = 0.5826986640974212

$ws.Range("B13").Value = 1.95795058464148
$ws.Range("C13").Value = 0.7302763736279871
$ws.Range("D13").Value = 0.1182647819889837
$ws.Range("E13").Value = 0.6989768860392616
$ws.Range("G13").Value = 1.065372199042287
$ws.Range("H13").Value = 0.8550754199502819
$ws.Range("I13").Value = 0.5820210209279111

$ws.Range("B14").Value = 1.925717962338638
$ws.Range("C14").Value = 0.7176348051098103
$ws.Range("D14").Value = 0.1165307235312127
$ws.Range("E14").Value = 0.6866713453964763
$ws.Range("G14").Value = 1.054655314732884
$ws.Range("H14").Value = 0.8513093633710014
$ws.Range("I14").Value = 0.5798275207083563

$ws.Range("B15").Value = 1.905979168566148
$ws.Range("C15").Value = 0.7098928477166169
$ws.Range("D15").Value = 0.1154696144489691
$ws.Range("E15").Value = 0.6791394260746983
$ws.Range("G15").Value = 1.048111838976951
$ws.Range("H15").Value = 0.8490174237211647
$ws.Range("I15").Value = 0.5784965380328586

$ws.Range("B16").Value = 1.792933712186141
$ws.Range("C16").Value = 0.6655467876013859
$ws.Range("D16").Value = 0.1094047635699837
$ws.Range("E16").Value = 0.6360598042797818
$ws.Range("G16").Value = 1.010931231222287
$ws.Range("H16").Value = 0.8361098395073441
$ws.Range("I16").Value = 0.5710609180318826

$ws.Range("B17").Value = 1.723644813983412
$ws.Range("C17").Value = 0.6383586528295382
$ws.Range("D17").Value = 0.1056981947430415
$ws.Range("E17").Value = 0.6097027264817143
$ws.Range("G17").Value = 0.9884017653685646
$ws.Range("H17").Value = 0.828392075553495
$ws.Range("I17").Value = 0.5666689366248647

$ws.Range("B18").Value = 1.683813533922148
$ws.Range("C18").Value = 0.6227264838037172
$ws.Range("D18").Value = 0.1035713484634613
$ws.Range("E18").Value = 0.5945677078392038
$ws.Range("G18").Value = 0.975544758203057
$ws.Range("H18").Value = 0.8240261285359338
$ws.Range("I18").Value = 0.5642044830309345

$ws.Range("B19").Value = 1.670331055943791
$ws.Range("C19").Value = 0.6174346438580756
$ws.Range("D19").Value = 0.1028520989463715
$ws.Range("E19").Value = 0.5894474127494647
$ws.Range("G19").Value = 0.9712088687973335
$ws.Range("H19").Value = 0.8225603914360988
$ws.Range("I19").Value = 0.5633806010073172

$ws.Range("B20").Value = 1.73101846196505
$ws.Range("C20").Value = 0.641252273991654
$ws.Range("D20").Value = 0.1060922384139076
$ws.Range("E20").Value = 0.6125058816425764
$ws.Range("G20").Value = 0.9907895447019541
$ws.Range("H20").Value = 0.8292060644940875
$ws.Range("I20").Value = 0.5671300714640353

$ws.Range("B21").Value = 1.935184499547063
$ws.Range("C21").Value = 0.7213476526380873
$ws.Range("D21").Value = 0.11703983889808
$ws.Range("E21").Value = 0.6902846115565211
$ws.Range("G21").Value = 1.057798754640828
$ws.Range("H21").Value = 0.852412441263624
$ws.Range("I21").Value = 0.5804691757670355

$ws.Range("B22").Value = 2.068824159933172
$ws.Range("C22").Value = 0.773754708855563
$ws.Range("D22").Value = 0.1242415922959168
$ws.Range("E22").Value = 0.7413640747868016
$ws.Range("G22").Value = 1.102525888616782
$ws.Range("H22").Value = 0.8682434122472671
$ws.Range("I22").Value = 0.589749348446702

$ws.Range("B23").Value = 1.997478833510229
$ws.Range("C23").Value = 0.7457781294085635
$ws.Range("D23").Value = 0.1203935106531446
$ws.Range("E23").Value = 0.714078215171952
$ws.Range("G23").Value = 1.078567352738133
$ws.Range("H23").Value = 0.8597326977865976
$ws.Range("I23").Value = 0.584744267350402

$ws.Range("B24").Value = 1.72768482389688
$ws.Range("C24").Value = 0.6399440721629617
$ws.Range("D24").Value = 0.1059140784885813
$ws.Range("E24").Value = 0.6112385193526961
$ws.Range("G24").Value = 0.9897097325915922
$ws.Range("H24").Value = 0.8288378387846649
$ws.Range("I24").Value = 0.5669214039302375

$ws.Range("B25").Value = 1.437497525831304
$ws.Range("C25").Value = 0.5259984073898067
$ws.Range("D25").Value = 0.09048864546777224
$ws.Range("E25").Value = 0.5012456823349396
$ws.Range("G25").Value = 0.8977191713597676
$ws.Range("H25").Value = 0.7982999516459017
$ws.Range("I25").Value = 0.5500478280253702
